# Adjusting to urbs excel-format: rename and add missing columns.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Commodity sheet: insert a new "Type" column (all rows are "Stock")
# between the existing "Commodity" and "price" columns.
# ---------------------------------------------------------------------
$wsCommodity = $wb.Worksheets.Item("Commodity")
$wsCommodity.Columns.Item(3).Insert()
$wsCommodity.Range("C1").Value = "Type"
$wsCommodity.Range("C2").Value = "Stock"
$wsCommodity.Range("C3").Value = "Stock"
$wsCommodity.Range("C4").Value = "Stock"

# ---------------------------------------------------------------------
# Demand sheet: rename the "Left"/"Right" demand columns to be explicit
# about the commodity they refer to.
# ---------------------------------------------------------------------
$wsDemand = $wb.Worksheets.Item("Demand")
$wsDemand.Range("B1").Value = "Left.Elec"
$wsDemand.Range("C1").Value = "Right.Elec"

# ---------------------------------------------------------------------
# Restore/update the per-sheet selection & active sheet so the workbook
# opens back on the Commodity sheet (matching the new canonical layout).
# ---------------------------------------------------------------------
$wsDemand.Activate() | Out-Null
$wsDemand.Range("B2").Select() | Out-Null

$wsCommodity.Activate() | Out-Null
$wsCommodity.Range("C5").Select() | Out-Null
